$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.119.68'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '3.757.54'
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.24%  '

$ws.Range("D7").Value = '3.754.90'
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.44%  '

$ws.Range("E10").Value = '  +4.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.71'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000248'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("D15").Value = '4.388.02'
$ws.Range("E15").Value = '  +0.23%  '

$ws.Range("D16").Value = '3.761.53'
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("D17").Value = '69.147.73'
$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.56%  '

$ws.Range("E20").Value = '  -0.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '490.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.727'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.47%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.08%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000148'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("E28").Value = '  -1.50%  '

$ws.Range("E30").Value = '  -0.63%  '

$ws.Range("E31").Value = '  +2.02%  '

$ws.Range("E32").Value = '  -4.07%  '

$ws.Range("E33").Value = '  -0.99%  '

$ws.Range("D34").Value = '3.905.03'
$ws.Range("E34").Value = '  +0.12%  '

$ws.Range("D35").Value = '3.699.74'
$ws.Range("E35").Value = '  +0.28%  '

$ws.Range("E37").Value = '  +6.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.28%  '

$ws.Range("E39").Value = '  -0.40%  '

$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.11'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.00%  '

$ws.Range("E42").Value = '  +0.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '428.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.58'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.63%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '142.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.57%  '

$ws.Range("D51").Value = '2.800.42'
$ws.Range("E51").Value = '  -1.66%  '
